$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 150 ("Femacal de La Calera" /
# Zanahoria, Provincia de Quillota, fecha 44518). Every existing data row from
# the old row 150 down to the old row 237 shifts down by one (old row 237
# becomes the new row 238), so insert a fresh row at 150 first.
$ws.Rows.Item(150).Insert()

$ws.Range("A150").Value = 3
$ws.Range("B150").Value = "Femacal de La Calera"
$ws.Range("C150").Value = "Coquimbo"
$ws.Range("D150").Value2 = 44518
$ws.Range("E150").Value = 5
$ws.Range("F150").Value = 100114013
$ws.Range("G150").Value = "Zanahoria"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 360
$ws.Range("K150").Value = 7000
$ws.Range("L150").Value = 7500
$ws.Range("M150").Value = 7250
$ws.Range("N150").Value = '$/saco 20 kilos'
$ws.Range("O150").Value = "Provincia de Quillota"
$ws.Range("P150").Value = 362
$ws.Range("Q150").Value = 20
$ws.Range("R150").Value = "Hortaliza"
